$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("General Settings")

# SerializationPath (named range -> D8) now calculates its own path instead
# of holding a hardcoded, machine-specific string. Use FormulaArray so the
# cell's existing number format / quote-prefix style is left untouched
# (matches how the workbook was actually edited upstream).
$ws.Range("D8").FormulaArray = '=SUBSTITUTE(LEFT(CELL("filename",A1),FIND("[",CELL("filename",A1),1)-1),"\XLS\","\XML\")'
